$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the easting/northing coordinates on row 3 to whole numbers
$ws.Range("Q3").Value = 407096
$ws.Range("R3").Value = 6702657

# Clear the start-time (Z3) and end-time (AB3) cells entirely
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
